# Fix English word errors in the database diagram sheet:
#  - "siteName"  -> "placeName"
#  - "siteImage" -> "placeImage"
#  - "yer_tip_id ... (FK - Place description)" -> "place_id ... (FK - Place description)"
# Also update the current selection/view (no frozen topLeftCell, select B11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B4: yer_tip_id ... -> place_id ...
$ws.Range("B4").Value = "place_id    int(10)  (FK - Place description)"

# B5: siteName -> placeName (keep original spacing/padding)
$ws.Range("B5").Value = "placeName         nvarchar(100)                 "

# B7: siteImage -> placeImage
$ws.Range("B7").Value = "placeImage        ByteArray"

# Update selection to B11 and scroll so column A is visible again (no frozen/topLeft override)
$ws.Range("B11").Select()
